$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.726.18"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.629.11"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "214.08"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").Value = "0.0631"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "19.43"
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "1.855.19"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "1.621.32"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "0.555"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "0.0₃0758"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").Value = "62.98"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "25.745.84"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "4.43"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "191.39"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").Value = "9.90"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").Value = "6.25"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "1.81"
$ws.Range("E25").Value = "  +2.96%  "
$ws.Range("D26").Value = "142.34"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "15.45"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "3.21"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").Value = "0.901"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "1.132.42"
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").Value = "0.540"
$ws.Range("E39").Value = "  -2.33%  "
$ws.Range("D40").Value = "0.0154"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "2.53"
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.51"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "100.03"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").Value = "0.796"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").Value = "1.764.85"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("D48").Value = "55.29"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "2.32"
$ws.Range("E51").Value = "  -7.08%  "
